$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.027510331597459
$ws.Range("D2").Value = 1.028816560393554
$ws.Range("E2").Value = 1.036184467811834
$ws.Range("F2").Value = 1.043970251305802
$ws.Range("J2").Value = 1.032667803777056
$ws.Range("K2").Value = 1.031632276888562
$ws.Range("L2").Value = 1.038978935124797
$ws.Range("M2").Value = 1.046742618042305
$ws.Range("N2").Value = 1.034134310507094

$ws.Range("C3").Value = 1.028893043923552
$ws.Range("D3").Value = 1.030039968870525
$ws.Range("E3").Value = 1.037455394122875
$ws.Range("F3").Value = 1.045422046661948
$ws.Range("J3").Value = 1.033688507682871
$ws.Range("K3").Value = 1.032662580953609
$ws.Range("L3").Value = 1.04005818633081
$ws.Range("M3").Value = 1.048003884081038
$ws.Range("N3").Value = 1.03515646392951

$ws.Range("C4").Value = 1.029787481076281
$ws.Range("D4").Value = 1.030831580599396
$ws.Range("E4").Value = 1.038277790879478
$ws.Range("F4").Value = 1.046361769813065
$ws.Range("J4").Value = 1.034348301310005
$ws.Range("K4").Value = 1.033328686352143
$ws.Range("L4").Value = 1.040756008892513
$ws.Range("M4").Value = 1.048819810016611
$ws.Range("N4").Value = 1.035817194539274

$ws.Range("C5").Value = 1.030163442453252
$ws.Range("D5").Value = 1.03116437379318
$ws.Range("E5").Value = 1.038623535380957
$ws.Range("F5").Value = 1.046756909016545
$ws.Range("J5").Value = 1.034625521237687
$ws.Range("K5").Value = 1.033608583556462
$ws.Range("L5").Value = 1.041049250978025
$ws.Range("M5").Value = 1.049162781780916
$ws.Range("N5").Value = 1.036094808151051

$ws.Range("C6").Value = 1.030226564640103
$ws.Range("D6").Value = 1.031220251308486
$ws.Range("E6").Value = 1.038681588028801
$ws.Range("F6").Value = 1.046823259419736
$ws.Range("J6").Value = 1.034672058548405
$ws.Range("K6").Value = 1.03365557178718
$ws.Range("L6").Value = 1.041098480522787
$ws.Range("M6").Value = 1.049220365808535
$ws.Range("N6").Value = 1.036141411550087

$ws.Range("C7").Value = 1.029792504920549
$ws.Range("D7").Value = 1.030836027395723
$ws.Range("E7").Value = 1.038282410695054
$ws.Range("F7").Value = 1.046367049364352
$ws.Range("J7").Value = 1.034352006150541
$ws.Range("K7").Value = 1.033332426873143
$ws.Range("L7").Value = 1.040759927685063
$ws.Range("M7").Value = 1.048824392988328
$ws.Range("N7").Value = 1.035820904641108

$ws.Range("C8").Value = 1.027977684463047
$ws.Range("D8").Value = 1.029230022185555
$ws.Range("E8").Value = 1.036613980787314
$ws.Range("F8").Value = 1.044460829820965
$ws.Range("J8").Value = 1.033012896046097
$ws.Range("K8").Value = 1.03198059285211
$ws.Range("L8").Value = 1.039343783782489
$ws.Range("M8").Value = 1.047168911987314
$ws.Range("N8").Value = 1.034479892846718

$ws.Range("C9").Value = 1.02477746014151
$ws.Range("D9").Value = 1.026399770131655
$ws.Range("E9").Value = 1.03367400778006
$ws.Range("F9").Value = 1.041104043147697
$ws.Range("J9").Value = 1.030647942229689
$ws.Range("K9").Value = 1.029593986455952
$ws.Range("L9").Value = 1.036844191959806
$ws.Range("M9").Value = 1.044250060285979
$ws.Range("N9").Value = 1.032111580524626

$ws.Range("C10").Value = 1.022642151200195
$ws.Range("D10").Value = 1.024512548296932
$ws.Range("E10").Value = 1.031713812199103
$ws.Range("F10").Value = 1.038867404259758
$ws.Range("J10").Value = 1.029067567110578
$ws.Range("K10").Value = 1.027999698602441
$ws.Range("L10").Value = 1.035174795042702
$ws.Range("M10").Value = 1.042302789122539
$ws.Range("N10").Value = 1.030528961091575

$ws.Range("C11").Value = 1.021717040028188
$ws.Range("D11").Value = 1.023695221703369
$ws.Range("E11").Value = 1.030864925184077
$ws.Range("F11").Value = 1.037899143375257
$ws.Range("J11").Value = 1.028382316426072
$ws.Range("K11").Value = 1.027308548652399
$ws.Range("L11").Value = 1.034451170293382
$ws.Range("M11").Value = 1.041459220739039
$ws.Range("N11").Value = 1.029842737272497

$ws.Range("C12").Value = 1.021373330462032
$ws.Range("D12").Value = 1.023391604010616
$ws.Range("E12").Value = 1.0305495899387
$ws.Range("F12").Value = 1.037539515714689
$ws.Range("J12").Value = 1.028127639064292
$ws.Range("K12").Value = 1.027051698930343
$ws.Range("L12").Value = 1.034182264958065
$ws.Range("M12").Value = 1.041145819046869
$ws.Range("N12").Value = 1.029587698239652

$ws.Range("C13").Value = 1.021447061137352
$ws.Range("D13").Value = 1.0234567322883
$ws.Range("E13").Value = 1.030617231431484
$ws.Range("F13").Value = 1.037616655907344
$ws.Range("J13").Value = 1.028182274834975
$ws.Range("K13").Value = 1.027106799823219
$ws.Range("L13").Value = 1.034239951509002
$ws.Range("M13").Value = 1.041213047692358
$ws.Range("N13").Value = 1.029642411599396

$ws.Range("C14").Value = 1.021688630595919
$ws.Range("D14").Value = 1.023670125128552
$ws.Range("E14").Value = 1.030838859924248
$ws.Range("F14").Value = 1.03786941590941
$ws.Range("J14").Value = 1.028361267691008
$ws.Range("K14").Value = 1.027287319963415
$ws.Range("L14").Value = 1.034428944942706
$ws.Range("M14").Value = 1.041433316166309
$ws.Range("N14").Value = 1.029821658645816

$ws.Range("C15").Value = 1.021837458463471
$ws.Range("D15").Value = 1.02380159990479
$ws.Range("E15").Value = 1.030975409664376
$ws.Range("F15").Value = 1.038025153180162
$ws.Range("J15").Value = 1.028471531801369
$ws.Range("K15").Value = 1.027398527590114
$ws.Range("L15").Value = 1.034545374162417
$ws.Range("M15").Value = 1.041569022376539
$ws.Range("N15").Value = 1.02993207934386

$ws.Range("C16").Value = 1.022703536584787
$ws.Range("D16").Value = 1.024566788122119
$ws.Range("E16").Value = 1.031770147356644
$ws.Range("F16").Value = 1.038931668698935
$ws.Range("J16").Value = 1.029113024854575
$ws.Range("K16").Value = 1.028045550504569
$ws.Range("L16").Value = 1.035222803180625
$ws.Range("M16").Value = 1.042358765423203
$ws.Range("N16").Value = 1.030574483390782

$ws.Range("C17").Value = 1.023246664388997
$ws.Range("D17").Value = 1.025046728201882
$ws.Range("E17").Value = 1.032268632942376
$ws.Range("F17").Value = 1.039500356710895
$ws.Range("J17").Value = 1.029515162702888
$ws.Range("K17").Value = 1.02845119063878
$ws.Range("L17").Value = 1.035647528930457
$ws.Range("M17").Value = 1.042854043434037
$ws.Range("N17").Value = 1.030977192320966

$ws.Range("C18").Value = 1.023563412991789
$ws.Range("D18").Value = 1.025326655217996
$ws.Range("E18").Value = 1.032559380893348
$ws.Range("F18").Value = 1.039832084249133
$ws.Range("J18").Value = 1.029749632680744
$ws.Range("K18").Value = 1.028687715498335
$ws.Range("L18").Value = 1.035895190745699
$ws.Range("M18").Value = 1.043142893684339
$ws.Range("N18").Value = 1.031211995273084

$ws.Range("C19").Value = 1.02367140799919
$ws.Range("D19").Value = 1.025422100860331
$ws.Range("E19").Value = 1.03265851682293
$ws.Range("F19").Value = 1.039945198563258
$ws.Range("J19").Value = 1.029829565696024
$ws.Range("K19").Value = 1.028768351217203
$ws.Range("L19").Value = 1.035979624672345
$ws.Range("M19").Value = 1.043241378014307
$ws.Range("N19").Value = 1.031292041802415

$ws.Range("C20").Value = 1.023188396991427
$ws.Range("D20").Value = 1.024995236643481
$ws.Range("E20").Value = 1.032215151250651
$ws.Range("F20").Value = 1.039439339691132
$ws.Range("J20").Value = 1.029472026479032
$ws.Range("K20").Value = 1.02840767740486
$ws.Range("L20").Value = 1.035601967512463
$ws.Range("M20").Value = 1.042800908656557
$ws.Range("N20").Value = 1.030933994838725

$ws.Range("C21").Value = 1.021617496723347
$ws.Range("D21").Value = 1.023607287012213
$ws.Range("E21").Value = 1.030773596466803
$ws.Range("F21").Value = 1.037794983670888
$ws.Range("J21").Value = 1.028308562782891
$ws.Range("K21").Value = 1.027234164785899
$ws.Range("L21").Value = 1.034373294400772
$ws.Range("M21").Value = 1.041368454362952
$ws.Range("N21").Value = 1.029768878890683

$ws.Range("C22").Value = 1.020629331602536
$ws.Range("D22").Value = 1.022734474334226
$ws.Range("E22").Value = 1.029867111170913
$ws.Range("F22").Value = 1.036761268503874
$ws.Range("J22").Value = 1.027576207601986
$ws.Range("K22").Value = 1.026495600789228
$ws.Range("L22").Value = 1.033600089044928
$ws.Range("M22").Value = 1.040467447624598
$ws.Range("N22").Value = 1.02903548368142

$ws.Range("C23").Value = 1.021153223831285
$ws.Range("D23").Value = 1.023197184788924
$ws.Range("E23").Value = 1.03034766926263
$ws.Range("F23").Value = 1.037309247385818
$ws.Range("J23").Value = 1.02796452379171
$ws.Range("K23").Value = 1.02688719799077
$ws.Range("L23").Value = 1.034010046548583
$ws.Range("M23").Value = 1.040945124589173
$ws.Range("N23").Value = 1.029424351324673

$ws.Range("C24").Value = 1.02321472567028
$ws.Range("D24").Value = 1.025018503505078
$ws.Range("E24").Value = 1.032239317356518
$ws.Range("F24").Value = 1.039466910589366
$ws.Range("J24").Value = 1.029491518160383
$ws.Range("K24").Value = 1.028427339403589
$ws.Range("L24").Value = 1.035622554983062
$ws.Range("M24").Value = 1.042824918089506
$ws.Range("N24").Value = 1.030953514200498

$ws.Range("C25").Value = 1.025605095095596
$ws.Range("D25").Value = 1.027131511033063
$ws.Range("E25").Value = 1.034434082541934
$ws.Range("F25").Value = 1.041971618222793
$ws.Range("J25").Value = 1.031259984173877
$ws.Range("K25").Value = 1.030211534456443
$ws.Range("L25").Value = 1.037490910712362
$ws.Range("M25").Value = 1.045004877348884
$ws.Range("N25").Value = 1.032724491638578

